$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B54").Value = "Administrador"
$ws.Range("C54").Value = "quiero saber que productos fueron los que mas llamaron la atención del usuario "
$ws.Range("D54").Value = "para poder tener mayores ventas "
$ws.Range("E54").Value = "Jesica Amaya "

$ws.Range("A55").Value = "US0024"
$ws.Range("B55").Value = "Administrador"
$ws.Range("C55").Value = "necesito q en cada producto me informe cuantas unidades en stock quedan disponibles"
$ws.Range("D55").Value = "para poder hacer control y reposición de cada producto "
$ws.Range("E55").Value = "Jesica Amaya"
